$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 18 (sheet row 19): issue about StoreKeeper screen - now resolved, "Won't Fix" by Sabethan
$ws.Range("G19").Value = "Won't Fix"
$ws.Range("H19").Value = "Sabethan"

# Row 34 (sheet row 36): Test class for storekeeper added - Fixed by Sabethan
$ws.Range("G36").Value = "Fixed"
$ws.Range("H36").Value = "Sabethan"

# Row 35 (sheet row 37): Vendor name uniqueness - Fixed by Sabethan
$ws.Range("G37").Value = "Fixed"
$ws.Range("H37").Value = "Sabethan"

# Row 36 (sheet row 38): Add Product dialog centering - Fixed by Sabethan
$ws.Range("G38").Value = "Fixed"
$ws.Range("H38").Value = "Sabethan"

# Restore selection/scroll position to reflect the latest review location
$ws.Range("H20").Select()
